# Edit script for PROJECT_NSRR_report.docx
# Applies:
#   1. A "complex-script" (BiDi) language tag on the title-image run (lang bidi=ne-NP)
#   2. Merges the two "A course project submitted ... [CT 401" runs into one run
#   3. Corrects "Rishikesh Poudel" -> "Rishikesh Paudel" and relocates the
#      "_GoBack" bookmark (Word always keeps only a single _GoBack, parked at
#      the location of the most recent edit) from the empty paragraph near the
#      DATE line to right after "Rishikesh Pa" in the name.

$d = $word.ActiveDocument

# --- 1. Tag the title image run's run properties with the BiDi language ---
$titleImgRange = $d.Range(0, 1)
$titleImgRange.LanguageIDOther = "ne-NP"

# --- 2. Merge the two bold/underline runs that make up the subtitle line
#        into a single run (identical formatting, so Find/Replace folds them
#        into one <w:r>) ---
$oldSubtitle = "A course project submitted to the department of Electronics and Computer Engineering in partial fulfillment of the requirements for the practical course on Computer Programming [CT 401"
$newSubtitle = "A course project submitted to the department of Electronics and Computer Engineering in partial fulfillment of the requirements for the practical course on Computer Programming [CT 401"
$d.Content.Find.Execute($oldSubtitle, $false, $false, $false, $false, $false, $true, 1, $false, $newSubtitle, 2)

# --- 3. Fix the submitter's name "Poudel" -> "Paudel" ---
$d.Content.Find.Execute("Poudel", $false, $false, $false, $false, $false, $true, 1, $false, "Paudel", 2)

# Relocate the _GoBack bookmark to sit right after "Rishikesh Pa" (i.e. right
# where the "o" was changed to "a"). Adding a bookmark named "_GoBack" moves
# the single allowed instance of it here and drops it from its previous spot.
$nameRange = $d.Content
$nameRange.Find.Execute("Rishikesh Pa", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($nameRange.End, $nameRange.End)
$d.Bookmarks.Add("_GoBack", $splitPoint)
